# Update docs/epexspot_prices.xlsx:
#  - "Prix Spot": append a new column BS ("23-aug") with 24 price values
#  - "Gaz": append a new row 68 (2025-08-21, 31.975)
#  - "CO2": append a new row 68 (2025-08-21, 71.5)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": new column BS = "23-aug"
# ---------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Bring over the header cell's formatting (bold, border, centered -
# the same style already used by the rest of row 1) before writing
# the new header text.
$wsSpot.Range("BR1").Copy()
$wsSpot.Range("BS1").PasteSpecial(-4122)
$wsSpot.Range("BS1").Value = "23-aug"

$spotValues = @(
    95.33,
    89.31999999999999,
    80.94,
    78.20999999999999,
    76.06,
    79.09999999999999,
    84.28,
    86.42,
    79.15000000000001,
    55.02,
    15.85,
    3.46,
    1.82,
    0.75,
    0.65,
    3.52,
    11.99,
    38.38,
    76.90000000000001,
    106.96,
    118.85,
    117.12,
    111.16,
    105.3
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 71).Value = $spotValues[$i]
}

# ---------------------------------------------------------------
# Sheet "Gaz": new row 68 = 2025-08-21 / 31.975
# ---------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$gazDate = $wsGaz.Cells.Item(68, 1)
$gazDate.NumberFormat = "@"
$gazDate.Value = "2025-08-21"
$gazDate.Style = "Normal"

$wsGaz.Cells.Item(68, 2).Value = 31.975

# ---------------------------------------------------------------
# Sheet "CO2": new row 68 = 2025-08-21 / 71.5
# ---------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$co2Date = $wsCO2.Cells.Item(68, 1)
$co2Date.NumberFormat = "@"
$co2Date.Value = "2025-08-21"
$co2Date.Style = "Normal"

$wsCO2.Cells.Item(68, 2).Value = 71.5
